$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue 'D2' '30.213.58'
Set-TextValue 'E2' '  +0.13%  '
Set-TextValue 'D3' '1.870.17'
Set-TextValue 'E3' '  +0.89%  '
Set-TextValue 'E4' '  +0.05%  '
Set-TextValue 'D5' '234.66'
Set-TextValue 'E5' '  -0.29%  '
Set-TextValue 'E6' '  +0.03%  '
Set-TextValue 'D7' '0.4700'
Set-TextValue 'E7' '  +0.16%  '
Set-TextValue 'D8' '0.2846'
Set-TextValue 'E8' '  -1.45%  '
Set-TextValue 'B9' 'OKB'
Set-TextValue 'C9' 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 'D9' '41.46'
Set-TextValue 'E9' '  -2.87%  '
Set-TextValue 'B10' 'Dogecoin'
Set-TextValue 'C10' 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextValue 'D10' '0.06565'
Set-TextValue 'E10' '  +0.24%  '
Set-TextValue 'B11' 'Solana'
Set-TextValue 'C11' 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextValue 'D11' '21.25'
Set-TextValue 'E11' '  -2.74%  '
Set-TextValue 'B12' 'TRON'
Set-TextValue 'C12' 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue 'D12' '0.07781'
Set-TextValue 'E12' '  -2.44%  '
Set-TextValue 'B13' 'Litecoin'
Set-TextValue 'C13' 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue 'D13' '96.31'
Set-TextValue 'E13' '  -1.19%  '
Set-TextValue 'B14' 'WrappedEther'
Set-TextValue 'C14' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D14' '1.870.34'
Set-TextValue 'E14' '  +0.74%  '
Set-TextValue 'B15' 'Polygon'
Set-TextValue 'C15' 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue 'D15' '0.6901'
Set-TextValue 'E15' '  +2.01%  '
Set-TextValue 'B16' 'Polkadot'
Set-TextValue 'C16' 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 'D16' '5.093'
Set-TextValue 'E16' '  -0.22%  '
Set-TextValue 'B17' 'BitcoinCash'
Set-TextValue 'C17' 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue 'D17' '266.97'
Set-TextValue 'E17' '  -0.85%  '
Set-TextValue 'B18' 'WrappedBTC'
Set-TextValue 'C18' 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue 'D18' '30.203.57'
Set-TextValue 'E18' '  +0.18%  '
Set-TextValue 'B19' 'Avalanche'
Set-TextValue 'C19' 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue 'D19' '13.68'
Set-TextValue 'E19' '  +0.49%  '
Set-TextValue 'B20' 'ShibaInu'
Set-TextValue 'C20' 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue 'D20' '0.000007730'
Set-TextValue 'E20' '  +0.38%  '
Set-TextValue 'B21' 'Dai'
Set-TextValue 'C21' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 'D21' '0.9998'
Set-TextValue 'E21' '  -0.02%  '
Set-TextValue 'B22' 'WrappedliquidstakedEther2.0'
Set-TextValue 'C22' 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue 'D22' '2.115.03'
Set-TextValue 'E22' '  +0.95%  '
Set-TextValue 'B23' 'BinanceUSD'
Set-TextValue 'C23' 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue 'D23' '0.9998'
Set-TextValue 'E23' '  -0.02%  '
Set-TextValue 'B24' 'Uniswap'
Set-TextValue 'C24' 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue 'D24' '5.229'
Set-TextValue 'E24' '  +0.50%  '
Set-TextValue 'B25' 'Chainlink'
Set-TextValue 'C25' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D25' '6.157'
Set-TextValue 'E25' '  +0.32%  '
Set-TextValue 'B26' 'Cosmos'
Set-TextValue 'C26' 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue 'D26' '9.469'
Set-TextValue 'E26' '  +3.38%  '
Set-TextValue 'B27' 'Monero'
Set-TextValue 'C27' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D27' '165.83'
Set-TextValue 'E27' '  -0.80%  '
Set-TextValue 'B28' 'EthereumClassic'
Set-TextValue 'C28' 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 'D28' '18.73'
Set-TextValue 'E28' '  -0.66%  '
Set-TextValue 'B29' 'LidoDAOToken'
Set-TextValue 'C29' 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue 'D29' '1.937'
Set-TextValue 'E29' '  +0.31%  '
Set-TextValue 'B30' 'Toncoin'
Set-TextValue 'C30' 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 'D30' '1.370'
Set-TextValue 'E30' '  -0.50%  '
Set-TextValue 'B31' 'Stellar'
Set-TextValue 'C31' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 'D31' '0.09935'
Set-TextValue 'E31' '  +0.91%  '
Set-TextValue 'B32' 'Filecoin'
Set-TextValue 'C32' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D32' '4.359'
Set-TextValue 'E32' '  +1.76%  '
Set-TextValue 'B33' 'PancakeSwap'
Set-TextValue 'C33' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue 'D33' '1.457'
Set-TextValue 'E33' '  -0.37%  '
Set-TextValue 'B34' 'InternetComputer(DFINITY)'
Set-TextValue 'C34' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D34' '4.043'
Set-TextValue 'E34' '  +1.22%  '
Set-TextValue 'B35' 'Hedera'
Set-TextValue 'C35' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D35' '0.04727'
Set-TextValue 'E35' '  +0.68%  '
Set-TextValue 'B36' 'ARBITRUM'
Set-TextValue 'C36' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 'D36' '1.130'
Set-TextValue 'E36' '  +1.07%  '
Set-TextValue 'B37' 'ImmutableX'
Set-TextValue 'C37' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D37' '0.6996'
Set-TextValue 'E37' '  +0.10%  '
Set-TextValue 'B38' 'HuobiToken'
Set-TextValue 'C38' 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue 'D38' '2.718'
Set-TextValue 'E38' '  +0.44%  '
Set-TextValue 'B39' 'VeChain'
Set-TextValue 'C39' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D39' '0.01863'
Set-TextValue 'E39' '  -0.37%  '
Set-TextValue 'B40' 'MXToken'
Set-TextValue 'C40' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 'D40' '2.778'
Set-TextValue 'E40' '  +6.79%  '
Set-TextValue 'B41' 'FraxShare'
Set-TextValue 'C41' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D41' '6.232'
Set-TextValue 'E41' '  -1.47%  '
Set-TextValue 'B42' 'Aave'
Set-TextValue 'C42' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D42' '72.63'
Set-TextValue 'E42' '  -0.60%  '
Set-TextValue 'B43' 'RenderToken'
Set-TextValue 'C43' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D43' '1.935'
Set-TextValue 'E43' '  +0.11%  '
Set-TextValue 'B44' 'PaxDollar'
Set-TextValue 'C44' 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue 'D44' '1.000'
Set-TextValue 'E44' '  +0.11%  '
Set-TextValue 'B45' 'TheSandbox'
Set-TextValue 'C45' 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue 'D45' '0.4144'
Set-TextValue 'E45' '  +0.35%  '
Set-TextValue 'B46' 'TrustWalletToken'
Set-TextValue 'C46' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 'D46' '0.8334'
Set-TextValue 'E46' '  -0.65%  '
Set-TextValue 'B47' 'Quant'
Set-TextValue 'C47' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue 'D47' '102.73'
Set-TextValue 'E47' '  -0.40%  '
Set-TextValue 'B48' 'Maker'
Set-TextValue 'C48' 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue 'D48' '970.32'
Set-TextValue 'E48' '  +3.52%  '
Set-TextValue 'B49' 'Aptos'
Set-TextValue 'C49' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 'D49' '7.068'
Set-TextValue 'E49' '  +0.71%  '
Set-TextValue 'B50' 'EnergySwap'
Set-TextValue 'C50' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D50' '9.187'
Set-TextValue 'E50' '  +0.13%  '
Set-TextValue 'B51' 'Elrond'
Set-TextValue 'C51' 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-TextValue 'D51' '34.47'
Set-TextValue 'E51' '  +1.81%  '
